# Update DB decimal types to int types on the "DB" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")

# Rows whose TYPE (column B) is "DECIMAL(10,0)" -> "INT(11)"
$decimal10Rows = @(5,22,26,27,38,39,40,44,56,59,60,61,62,63,64,65,66,77,78,79,80,92,94)
foreach ($r in $decimal10Rows) {
    $ws.Cells.Item($r, 2).Value = "INT(11)"
}

# Rows whose TYPE (column B) is "DECIMAL(1,0)" -> "INT(1)"
$decimal1Rows = @(16,32,50,71,86,99)
foreach ($r in $decimal1Rows) {
    $ws.Cells.Item($r, 2).Value = "INT(1)"
}

# Rows whose CONSTRAINTS (column C) is "NOT NULL PRIMARY KEY" -> "NOT NULL PRIMARY KEY AUTO_INCREMENT"
$primaryKeyRows = @(5,22,38,56,77,92)
foreach ($r in $primaryKeyRows) {
    $ws.Cells.Item($r, 3).Value = "NOT NULL PRIMARY KEY AUTO_INCREMENT"
}

# Column C width change (grew wider to fit "NOT NULL PRIMARY KEY AUTO_INCREMENT")
$ws.Columns.Item(3).ColumnWidth = 38.75

# Sheet view changes: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 63
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E90:E99").Select()

$wb.Save()
